$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.128.18"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "1.643.61"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "216.88"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "0.255"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").Value = "0.0626"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "19.96"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "1.873.32"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "1.641.56"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "0.544"
$ws.Range("E15").Value = "  +3.06%  "
$ws.Range("D16").Value = "67.27"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").Value = "27.119.23"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").Value = "0.0₃0738"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").Value = "218.38"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "6.82"
$ws.Range("E21").Value = "  +2.64%  "
$ws.Range("E22").Value = "  +5.76%  "
$ws.Range("D23").Value = "4.40"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("D24").Value = "9.18"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "147.87"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "7.52"
$ws.Range("E27").Value = "  +1.95%  "
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("D29").Value = "15.75"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").Value = "0.0507"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").Value = "1.260.42"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("E37").Value = "  +2.05%  "
$ws.Range("D38").Value = "0.543"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("E42").Value = "  +6.70%  "
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "1.783.35"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.28"
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("D45").Value = "61.71"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("D46").Value = "91.74"
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("E47").Value = "  +0.73%  "
$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.63"
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.0972"
$ws.Range("E51").Value = "  +0.12%  "
